# Auto commit 13-06-2025 14:36:39.85
#
# This script:
#  1. Fixes three mismatched item descriptions (wrong product codes that were
#     duplicated in the shared-string table) on the "MWSS with Compressor Pump"
#     and "Compressor Pump Installation" sheets.
#  2. Updates the saved selection (active cell) on a few sheets and switches
#     the active / tab-selected sheet from "200 mm Tubewell Construction" to
#     "MWSS with Compressor Pump".

$wb = $excel.ActiveWorkbook

$wsSubmersiblePump   = $wb.Worksheets.Item(1)   # MWSS with Submersible Pump
$wsCompressorPump    = $wb.Worksheets.Item(2)   # MWSS with Compressor Pump
$wsTubewell200       = $wb.Worksheets.Item(6)   # 200 mm Tubewell Construction
$wsCompressorInstall = $wb.Worksheets.Item(8)   # Compressor Pump Installation

# --- Correct the mismatched item text / product codes -----------------------

$wsCompressorPump.Range("A6").Value = "20 mm dia HDPE pipe (DG) (8kg) (Code: GWDMR076)"
$wsCompressorPump.Range("A7").Value = "32 mm dia UPVC pipe (Code: GWDMR072)"
$wsCompressorInstall.Range("A3").Value = "Compressor pump 2 HP single phase (Code: GWDMR065)"

# --- Update saved selections -------------------------------------------------

$wsSubmersiblePump.Range("A8").Select()
$wsCompressorInstall.Range("A6").Select()

# Select this one last so it ends up the active / tab-selected sheet.
$wsCompressorPump.Range("A12").Select()
